$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Transcriptions")
$ws2 = $wb.Worksheets.Item("Annotations")

# Remove the stray "trac" placeholder in Transcriptions!B15 first so its
# shared-string slot is freed up (and gets reclaimed by the next new string).
$ws1.Range("B15").Clear()

# --- Annotations: new reference rows 134-143 -----------------------------
# Row 134 - Deptford (Place)
$ws2.Range("A134").Value = "Deptford"
$ws2.Range("B134").Value = "Place"
$ws2.Range("C134").Value = "pla-dept"
$ws2.Range("D134").Value = "../resources/annotations.xml#pla-dept"
$ws2.Rows.Item(134).RowHeight = 32

# Row 135 - Frizer, Ingram (Person)
$ws2.Range("A135").Value = "Frizer, Ingram"
$ws2.Range("B135").Value = "Person"
$ws2.Range("C135").Value = "psn-ifriz"
$ws2.Range("D135").Value = "../resources/annotations.xml#psn-ifriz"
$ws2.Rows.Item(135).RowHeight = 32

# Row 136 - Jonson, Ben (Person)
$ws2.Range("A136").Value = "Jonson, Ben"
$ws2.Range("B136").Value = "Person"
$ws2.Range("C136").Value = "psn-bjon"
$ws2.Range("D136").Value = "../resources/annotations.xml#psn-bjon"
$ws2.Rows.Item(136).RowHeight = 32

# Row 137 - Spenser, Gabriel (Person)
$ws2.Range("A137").Value = "Spenser, Gabriel"
$ws2.Range("B137").Value = "Person"
$ws2.Range("C137").Value = "psn-gspe"
$ws2.Range("D137").Value = "../resources/annotations.xml#psn-gspe"
$ws2.Rows.Item(137).RowHeight = 32

# Row 138 - Kempe, William (Person)
$ws2.Range("A138").Value = "Kempe, William"
$ws2.Range("B138").Value = "Person"
$ws2.Range("C138").Value = "psn-wkem"
$ws2.Range("D138").Value = "../resources/annotations.xml#psn-wkem"
$ws2.Rows.Item(138).RowHeight = 32

# Row 139 - Norwich (Place)
$ws2.Range("A139").Value = "Norwich"
$ws2.Range("B139").Value = "Place"
$ws2.Range("C139").Value = "pla-norw"
$ws2.Range("D139").Value = "../resources/annotations.xml#pla-norw"
$ws2.Rows.Item(139).RowHeight = 32

# Rows 140-141 - Tennyson, Arthur / Browning, Robert (Person) - names then
# xml:ids then links were entered column-by-column across the pair.
$ws2.Range("A140").Value = "Tennyson, Arthur"
$ws2.Range("A141").Value = "Browning, Robert"
$ws2.Range("B140").Value = "Person"
$ws2.Range("B141").Value = "Person"
$ws2.Range("C140").Value = "psn-aten"
$ws2.Range("C141").Value = "psn-rbro"
$ws2.Range("D140").Value = "../resources/annotations.xml#psn-aten"
$ws2.Range("D141").Value = "../resources/annotations.xml#psn-rbro"
$ws2.Rows.Item(140).RowHeight = 32
$ws2.Rows.Item(141).RowHeight = 32

# Row 142 - Hobbes, Thomas (Person), with a note
$ws2.Range("A142").Value = "Hobbes, Thomas"
$ws2.Range("B142").Value = "Person"
$ws2.Range("C142").Value = "psn-thob"
$ws2.Range("D142").Value = "../resources/annotations.xml#psn-thob"
$ws2.Range("E142").Value = "Thoby?"
$ws2.Rows.Item(142).RowHeight = 32

# Row 143 - Garden, The (Literary Work), with a note
$ws2.Range("D143").Value = "../resources/annotations.xml#lit-gard"
$ws2.Range("C143").Value = "lit-gard"
$ws2.Range("A143").Value = "Garden, The"
$ws2.Range("B143").Value = "Literary Work"
$ws2.Range("E143").Value = "Marvell"
$ws2.Rows.Item(143).RowHeight = 32

# --- View state -------------------------------------------------------
$ws2.Range("E143").Select()
$ws2.Activate()
